# Auto-generated Excel COM-interop edit script
# Applies numeric updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4464.2905
$ws.Range("I28").Value = 444.44446
$ws.Range("K28").Value = 444.44446
$ws.Range("M28").Value = 40.55554000000001

$ws.Range("H32").Value = 13387.667
$ws.Range("J32").Value = 14436.125
$ws.Range("L32").Value = 14436.125
$ws.Range("N32").Value = -15088.125

$ws.Range("H51").Value = 14114.866
$ws.Range("I51").Value = 13385
$ws.Range("J51").Value = 14167
$ws.Range("K51").Value = 13385
$ws.Range("L51").Value = 14167
$ws.Range("M51").Value = -12901
$ws.Range("N51").Value = -15135

$ws.Range("H98").Value = 2061.8965
$ws.Range("I98").Value = 2111.3215
$ws.Range("K98").Value = 2111.3215
$ws.Range("M98").Value = -613.3215

$ws.Range("H106").Value = 2081.7693
$ws.Range("I106").Value = 1609
$ws.Range("K106").Value = 1609
$ws.Range("M106").Value = -978

$ws.Range("H112").Value = 2924.7568
$ws.Range("J112").Value = 2517.8
$ws.Range("L112").Value = 7553.400000000001
$ws.Range("N112").Value = -9769.400000000001

$ws.Range("H122").Value = 2061.8965
$ws.Range("I122").Value = 2111.3215
$ws.Range("K122").Value = 6333.9645
$ws.Range("M122").Value = -3883.9645

$ws.Range("H132").Value = 303134.6
$ws.Range("I132").Value = 360197.16
$ws.Range("K132").Value = 1080591.48
$ws.Range("M132").Value = -1078061.48

$ws.Range("H138").Value = 2718.5823
$ws.Range("I138").Value = 835.40424
$ws.Range("J138").Value = 5484.5
$ws.Range("K138").Value = 2506.21272
$ws.Range("L138").Value = 16453.5
$ws.Range("M138").Value = 2633.78728
$ws.Range("N138").Value = -26733.5

$ws.Range("H141").Value = 1021.9667
$ws.Range("I141").Value = 1060.5769
$ws.Range("J141").Value = 771
$ws.Range("K141").Value = 3181.7307
$ws.Range("L141").Value = 2313
$ws.Range("M141").Value = 1998.2693
$ws.Range("N141").Value = -12673

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4320.9
$ws.Range("I2").Value = 1561.091
$ws.Range("K2").Value = 1561.091
$ws.Range("M2").Value = -1448.091

$ws.Range("H61").Value = 5889.857
$ws.Range("I61").Value = 5642.3335
$ws.Range("K61").Value = 5642.3335
$ws.Range("M61").Value = -5430.3335

$ws.Range("H116").Value = 4320.9
$ws.Range("I116").Value = 1561.091
$ws.Range("K116").Value = 1561.091
$ws.Range("M116").Value = 732.9090000000001

$ws.Range("H122").Value = 2872.7827
$ws.Range("I122").Value = 2210.8
$ws.Range("J122").Value = 4114
$ws.Range("K122").Value = 6632.400000000001
$ws.Range("L122").Value = 12342
$ws.Range("M122").Value = -4182.400000000001
$ws.Range("N122").Value = -17242

$ws.Range("H132").Value = 451019.38
$ws.Range("I132").Value = 603649.4
$ws.Range("J132").Value = 54181.35
$ws.Range("K132").Value = 1810948.2
$ws.Range("L132").Value = 162544.05
$ws.Range("M132").Value = -1808418.2
$ws.Range("N132").Value = -167604.05

$ws.Range("H136").Value = 5889.857
$ws.Range("I136").Value = 5642.3335
$ws.Range("K136").Value = 16927.0005
$ws.Range("M136").Value = -14377.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4320.9
$ws.Range("I3").Value = 1561.091
$ws.Range("K3").Value = 1561.091
$ws.Range("M3").Value = -1447.091

$ws.Range("H22").Value = 20255.4
$ws.Range("I22").Value = 319.25
$ws.Range("J22").Value = 100000
$ws.Range("K22").Value = 319.25
$ws.Range("L22").Value = 100000
$ws.Range("M22").Value = -146.25
$ws.Range("N22").Value = -100346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H31").Value = 5964.391
$ws.Range("I31").Value = 3417.2
$ws.Range("K31").Value = 3417.2
$ws.Range("M31").Value = -3122.2

$ws.Range("H34").Value = 5964.391
$ws.Range("I34").Value = 3417.2
$ws.Range("K34").Value = 3417.2
$ws.Range("M34").Value = -3215.2

$ws.Range("H58").Value = 16671422
$ws.Range("I58").Value = 21279836
$ws.Range("K58").Value = 21279836
$ws.Range("M58").Value = -21279633

$ws.Range("H132").Value = 5406.095
$ws.Range("I132").Value = 4307.943
$ws.Range("J132").Value = 10896.857
$ws.Range("K132").Value = 12923.829
$ws.Range("L132").Value = 32690.571
$ws.Range("M132").Value = -10393.829
$ws.Range("N132").Value = -37750.571

$ws.Range("H134").Value = 55565556
$ws.Range("I134").Value = 111117960
$ws.Range("K134").Value = 333353880
$ws.Range("M134").Value = -333351345

$ws.Range("H136").Value = 16671422
$ws.Range("I136").Value = 21279836
$ws.Range("K136").Value = 63839508
$ws.Range("M136").Value = -63836958

$ws.Range("H141").Value = 74998.5
$ws.Range("J141").Value = 74998.5
$ws.Range("L141").Value = 74998.5
$ws.Range("N141").Value = -85358.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 409
$ws.Range("I60").Value = 79
$ws.Range("K60").Value = 237
$ws.Range("M60").Value = 14

$ws.Range("H68").Value = 836499.3
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 5997
$ws.Range("M68").Value = -5186

$ws.Range("H71").Value = 836499.3
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 17991
$ws.Range("M71").Value = -13935

$ws.Range("H113").Value = 3462565.5
$ws.Range("I113").Value = 13333767
$ws.Range("J113").Value = 501205
$ws.Range("K113").Value = 40001301
$ws.Range("L113").Value = 1503615
$ws.Range("M113").Value = -39999131
$ws.Range("N113").Value = -1507955

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9234.223
$ws.Range("I70").Value = 5393.2144
$ws.Range("K70").Value = 5393.2144
$ws.Range("M70").Value = -5123.2144

$ws.Range("H73").Value = 9234.223
$ws.Range("I73").Value = 5393.2144
$ws.Range("K73").Value = 5393.2144
$ws.Range("M73").Value = -4457.2144

$ws.Range("H97").Value = 972.2857
$ws.Range("I97").Value = 922.2593000000001
$ws.Range("J97").Value = 1141.125
$ws.Range("K97").Value = 922.2593000000001
$ws.Range("L97").Value = 1141.125
$ws.Range("M97").Value = -426.2593000000001
$ws.Range("N97").Value = -2133.125

$ws.Range("H102").Value = 3983.106
$ws.Range("I102").Value = 2742.4
$ws.Range("K102").Value = 2742.4
$ws.Range("M102").Value = -1120.4

$ws.Range("H107").Value = 1096.25
$ws.Range("J107").Value = 1159.4
$ws.Range("L107").Value = 1159.4
$ws.Range("N107").Value = -4999.4

$ws.Range("H141").Value = 43336.855
$ws.Range("J141").Value = 43336.855
$ws.Range("L141").Value = 43336.855
$ws.Range("N141").Value = -53696.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 788.3200000000001
$ws.Range("J22").Value = 986.1818
$ws.Range("L22").Value = 986.1818
$ws.Range("N22").Value = -1576.1818

$ws.Range("H27").Value = 788.3200000000001
$ws.Range("J27").Value = 986.1818
$ws.Range("L27").Value = 986.1818
$ws.Range("N27").Value = -1200.1818

$ws.Range("H40").Value = 3681.59
$ws.Range("I40").Value = 3597.8
$ws.Range("K40").Value = 3597.8
$ws.Range("M40").Value = -3461.8

$ws.Range("H122").Value = 2759.1562
$ws.Range("I122").Value = 2559.077
$ws.Range("K122").Value = 7677.231000000001
$ws.Range("M122").Value = -5227.231000000001

$ws.Range("H136").Value = 30616726
$ws.Range("I136").Value = 17246186
$ws.Range("K136").Value = 51738558
$ws.Range("M136").Value = -51736008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 3850096
$ws.Range("J23").Value = 6997.25
$ws.Range("L23").Value = 6997.25
$ws.Range("N23").Value = -7455.25

$ws.Range("H107").Value = 16675433
$ws.Range("I107").Value = 22224234
$ws.Range("K107").Value = 66672702
$ws.Range("M107").Value = -66670782

$ws.Range("H132").Value = 9882.130999999999
$ws.Range("I132").Value = 10402.866
$ws.Range("J132").Value = 8905.75
$ws.Range("K132").Value = 31208.598
$ws.Range("L132").Value = 26717.25
$ws.Range("M132").Value = -28678.598
$ws.Range("N132").Value = -31777.25

$ws.Range("H136").Value = 9437905
$ws.Range("I136").Value = 11628876
$ws.Range("J136").Value = 16731.9
$ws.Range("K136").Value = 34886628
$ws.Range("L136").Value = 50195.7
$ws.Range("M136").Value = -34884078
$ws.Range("N136").Value = -55295.7
